$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data dropped the "ECs" target-cluster row entirely (it was row 2,
# i.e. the Calca-Calcr -> ECs pair). Deleting the whole row shifts the rows
# below it up by one and lets the shared-strings table drop the now-unused
# "ECs" string on save (matching uniqueCount 25 -> 24).
$ws.Rows.Item(2).Delete()

# After the shift:
#   row 2 is the old "FAPs" row (was row 3)
#   row 3 is the old "MuSCs" row (was row 4)
# Both rows get refreshed TPM-derived numbers; row 3 additionally keeps the
# K/L values it already carried (3 and 1) from the old row 4, unchanged.

# Row 2 (target cluster: FAPs) - refreshed values
$ws.Cells.Item(2, 7).Value = 0.41903
$ws.Cells.Item(2, 8).Value = 1.25709
$ws.Cells.Item(2, 13).Value = 0.0005823333333333334
$ws.Cells.Item(2, 14).Value = 0.001747
$ws.Cells.Item(2, 15).Value = 0.00009625315715314126
$ws.Cells.Item(2, 16).Value = 0.00009625315715314125
$ws.Cells.Item(2, 17).Value = 0.0002440151366666667
$ws.Cells.Item(2, 18).Value = 0.00219613623
$ws.Cells.Item(2, 19).Value = 0.00009625315715314126
$ws.Cells.Item(2, 20).Value = 0.00009625315715314125

# Row 3 (target cluster: MuSCs) - refreshed values
$ws.Cells.Item(3, 7).Value = 0.41903
$ws.Cells.Item(3, 8).Value = 1.25709
$ws.Cells.Item(3, 13).Value = 6.049435666666667
$ws.Cells.Item(3, 14).Value = 18.148307
$ws.Cells.Item(3, 15).Value = 0.9999037468428469
$ws.Cells.Item(3, 16).Value = 0.9999037468428468
$ws.Cells.Item(3, 17).Value = 2.534895027403333
$ws.Cells.Item(3, 18).Value = 22.81405524662999
$ws.Cells.Item(3, 19).Value = 0.9999037468428469
$ws.Cells.Item(3, 20).Value = 0.9999037468428468
